$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 167-168 (existing rows 167-174 shift down to 169-176)
$ws.Rows("167:168").Insert()

# New row 167: Bing cherries, Primera quality
$ws.Range("A167").Value = 5
$ws.Range("B167").Value = "Macroferia Regional de Talca"
$ws.Range("C167").Value = "Maule"
$ws.Range("D167").Value = 44578
$ws.Range("E167").Value = 7
$ws.Range("F167").Value = "Fruta"
$ws.Range("G167").Value = 100103
$ws.Range("H167").Value = "Frutos de hueso (carozo)"
$ws.Range("I167").Value = 100103001
$ws.Range("J167").Value = "Cereza"
$ws.Range("K167").Value = "Bing"
$ws.Range("L167").Value = "Primera"
$ws.Range("M167").Value = 250
$ws.Range("N167").Value = 4500
$ws.Range("O167").Value = 5000
$ws.Range("P167").Value = 4800
$ws.Range("Q167").Value = "$/bandeja 10 kilos"
$ws.Range("R167").Value = "Provincia de Curicó"
$ws.Range("S167").Value = 480
$ws.Range("T167").Value = 10

# New row 168: Lapins cherries, Primera quality
$ws.Range("A168").Value = 5
$ws.Range("B168").Value = "Macroferia Regional de Talca"
$ws.Range("C168").Value = "Maule"
$ws.Range("D168").Value = 44578
$ws.Range("E168").Value = 7
$ws.Range("F168").Value = "Fruta"
$ws.Range("G168").Value = 100103
$ws.Range("H168").Value = "Frutos de hueso (carozo)"
$ws.Range("I168").Value = 100103001
$ws.Range("J168").Value = "Cereza"
$ws.Range("K168").Value = "Lapins"
$ws.Range("L168").Value = "Primera"
$ws.Range("M168").Value = 150
$ws.Range("N168").Value = 5000
$ws.Range("O168").Value = 5000
$ws.Range("P168").Value = 5000
$ws.Range("Q168").Value = "$/bandeja 10 kilos"
$ws.Range("R168").Value = "Provincia de Curicó"
$ws.Range("S168").Value = 500
$ws.Range("T168").Value = 10
